# Username/password changed from MHRA12345 TO MHRA123456
# - Update the password on the "Auto.*"/"Noor.Uddin.*" rows (B2:B7) on Sheet1
#   from "MHRA12345" to "MHRA123456".
# - Sheet1 becomes the active tab / sheet (selection moves to D13).
# - Column B on Sheet1 widens by one character (10 -> 11) to fit the longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the six password cells that held the old value.
$ws.Range("B2").Value = "MHRA123456"
$ws.Range("B3").Value = "MHRA123456"
$ws.Range("B4").Value = "MHRA123456"
$ws.Range("B5").Value = "MHRA123456"
$ws.Range("B6").Value = "MHRA123456"
$ws.Range("B7").Value = "MHRA123456"

# Widen column B to fit the new, longer password text.
$ws.Columns.Item(2).ColumnWidth = 10.166666666666666

# Make Sheet1 the active sheet/tab and move the selection to D13.
$ws.Activate()
$ws.Range("D13").Select()
